$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($r, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# Update the "last updated" timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 28 de Mayo de 2020 a las 16:10"

# Row 12 / Row 13 swapped rank (India moves above Turquia) with refreshed counts
Set-Row 12 @("India", 160310, 2224, 68713, 87037, 0, 26, 4560)
Set-Row 13 @("Turquia", 159797, 0, 122793, 32573, 0, 0, 4431)

# Row 41 refreshed counts (country unchanged)
$ws.Cells.Item(41, 5).Value = 4931
$ws.Cells.Item(41, 7).Value = 4
$ws.Cells.Item(41, 8).Value = 1231

# Row 66 refreshed counts (country unchanged)
$ws.Cells.Item(66, 4).Value = 6580
$ws.Cells.Item(66, 5).Value = 467

# Row 68 / Row 69 swapped rank (Irak moves above Camerun) with refreshed counts
Set-Row 68 @("Irak", 5457, 322, 2971, 2307, 0, 4, 179)
Set-Row 69 @("Camerun", 5436, 0, 1996, 3265, 0, 0, 175)

# Row 95 / Row 96 swapped rank (Mayotte moves above Lituania) with refreshed counts
Set-Row 95 @("Mayotte", 1670, 25, 1315, 334, 0, 1, 21)
Set-Row 96 @("Lituania", 1656, 9, 1193, 395, 0, 2, 68)

# Row 197 / Row 198 swapped rank (Fiyi moves above Curazao)
$ws.Cells.Item(197, 1).Value = "Fiyi"
$ws.Cells.Item(197, 4).Value = 15
$ws.Cells.Item(197, 8).Value = 0

$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 8).Value = 1

# Row 199 renamed (Santa Lucia -> Nueva Caledonia), rows 200/201 shift down accordingly
$ws.Cells.Item(199, 1).Value = "Nueva Caledonia"

$ws.Cells.Item(200, 1).Value = "Santa Lucia"
$ws.Cells.Item(200, 4).Value = 18
$ws.Cells.Item(200, 8).Value = 0

$ws.Cells.Item(201, 1).Value = "Belice"
$ws.Cells.Item(201, 4).Value = 16
$ws.Cells.Item(201, 8).Value = 2

# Row 210 / Row 211 swapped rank (Seychelles moves above Montserrat)
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1

# Row 213 / Row 214 swapped rank (Islas Virgenes Britanicas moves above Papua Nueva Guinea)
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0

# Row 215 / Row 216 renamed (San Bartolome <-> Bonaire, San Eustaquio y Saba)
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216, 1).Value = "San Bartolome"
